$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '49.688.84'
$ws.Range("E2").Value = '  -0.63%  '

$ws.Range("D3").Value = '2.652.60'
$ws.Range("E3").Value = '  +0.14%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '112.84'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.23%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '327.39'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.27%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.524'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.04%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.551'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.92%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.79'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.21%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.04'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.64%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0817'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.86%  '

$ws.Range("E13").Value = '  +2.17%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.61'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.96%  '

$ws.Range("D15").Value = '3.065.74'
$ws.Range("E15").Value = '  +0.17%  '

$ws.Range("D16").Value = '2.635.03'
$ws.Range("E16").Value = '  -1.41%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.860'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.51%  '

$ws.Range("D18").Value = '49.619.17'
$ws.Range("E18").Value = '  -0.56%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.34'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.88%  '

$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.70'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.34%  '

$ws.Range("B21").Value = 'ImmutableX'
$ws.Range("C21").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.91'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.31%  '

$ws.Range("D22").Value = '0.0₃0951'
$ws.Range("E22").Value = '  -0.75%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '269.01'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.66%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.12'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.14%  '

$ws.Range("E25").Value = '  -0.48%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.16'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.49%  '

$ws.Range("E27").Value = '  +0.05%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.19'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.53%  '

$ws.Range("E29").Value = '  -0.76%  '

$ws.Range("E30").Value = '  -2.17%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.91'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.66%  '

$ws.Range("E32").Value = '  -1.52%  '

$ws.Range("E33").Value = '  +0.14%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0822'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.55%  '

$ws.Range("E35").Value = '  -2.09%  '

$ws.Range("E36").Value = '  -0.23%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.94'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.60%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.04'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.61%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.13'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.87%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '129.43'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.97%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '23.77'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.93%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0343'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +8.61%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.29'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.63%  '

$ws.Range("E44").Value = '  -0.59%  '

$ws.Range("D45").Value = '2.065.99'
$ws.Range("E45").Value = '  -0.98%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.32'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.37%  '

$ws.Range("E47").Value = '  +5.79%  '

$ws.Range("E48").Value = '  -3.54%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.94'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.19%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.27'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.48%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '59.02'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.01%  '
